# Reformat code & add some missing 'Step' info
# The "Result" column (D) on the "data" sheet didn't have an outcome
# recorded for the "Only me" / "Public" privacy-type test rows yet.
# Fill those in as failing steps.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "Fail"
$ws.Range("D3").Value = "Fail"
